$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row values (D1 reuses "AWB", new headers appended in order) ---
$ws.Range("D1").Value = "AWB"
$ws.Range("E1").Value = "Manifest Pieces(Please Enter Numaric Value)"
$ws.Range("F1").Value = "Manifest Weight(Please Enter Numaric Value)"
$ws.Range("G1").Value = "Shipper’s Name "
$ws.Range("H1").Value = "Shipper’s Address "
$ws.Range("I1").Value = "Shipper’s Contact Number "
$ws.Range("J1").Value = "Consignee's Name"
$ws.Range("K1").Value = "Consignee's Address"
$ws.Range("L1").Value = "Consignee's Contact Number"
$ws.Range("M1").Value = "Origin"
$ws.Range("N1").Value = "Destination"
$ws.Range("O1").Value = "Cargo Type"
$ws.Range("P1").Value = "Advisory To Consignee"

# New header cells need the same wrap-text formatting as the existing D1:F1 headers
$ws.Range("G1:P1").WrapText = $true

# --- Column widths for the new columns (chosen so the stored/serialized
#     character width lands on the closest achievable value to the target) ---
$ws.Columns.Item(5).ColumnWidth = 20.16666666666667
$ws.Columns.Item(6).ColumnWidth = 19.16666666666667
$ws.Range("G1:H1").EntireColumn.ColumnWidth = 17.66666666666667
$ws.Columns.Item(9).ColumnWidth = 18
$ws.Columns.Item(10).ColumnWidth = 16.5
$ws.Columns.Item(11).ColumnWidth = 19.5
$ws.Columns.Item(12).ColumnWidth = 18.83333333333334
$ws.Columns.Item(13).ColumnWidth = 11.16666666666667
$ws.Columns.Item(14).ColumnWidth = 17.33333333333334
$ws.Columns.Item(15).ColumnWidth = 12.5
$ws.Columns.Item(16).ColumnWidth = 13.83333333333333

# --- Row height for the (now taller) header row ---
$ws.Rows.Item(1).RowHeight = 60

# --- View / selection state ---
$ws.Range("M18").Select()
$excel.ActiveWindow.ScrollColumn = 4
